$d = $word.ActiveDocument

# Update the title headline (appears in Heading1 and also later in a bold run)
$d.Content.Find.Execute(
    "Play Hot 4 Cash for Free - Retro Slot with Unique Bonus Features",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Play Hot 4 Cash for Free", 2)

# Update the "What we like" bullet points
$d.Content.Find.Execute(
    "Unique Hot Zones bonus feature",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Unique design and retro theme", 2)

$d.Content.Find.Execute(
    "3 bonus rounds with generous multipliers",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Decent payout and bonus features", 2)

$d.Content.Find.Execute(
    "Retro design with immersive visuals",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Hot Zones trigger bonus rounds", 2)

$d.Content.Find.Execute(
    "Lucrative payout values for symbols",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Wild symbol replaces other symbols", 2)

# Update the "What we don't like" bullet point
$d.Content.Find.Execute(
    "Limited number of paylines",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Only three bonus rounds", 2)

# Update the meta description line
$d.Content.Find.Execute(
    "Review of Hot 4 Cash - Play this retro-themed slot with 10 paylines, three bonus features, and Hot Zones that trigger bonus rounds, all for free.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Read our review of Hot 4 Cash and play this retro slot game for free.", 2)
